$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1.03
$ws.Range("C4").Value = "cleaning up the assets mod"
$ws.Range("B4").Value = "logger enhancement"

$ws.Range("A5").Value = 1.04
$ws.Range("B5").Value = "pledges setup"

$ws.Range("B6").Value = "arrangements setup"

$ws.Range("C9").Select()
